# "Generate Report for handoff"
#
# This applies the localization-status report regeneration:
#   - Status text for the source file changes from "Handoff transform failed"
#     to "Ready for handoff" on both the zh-cn and de-de sheets.
#   - A freshly generated handoff (.xlf) file is now associated with the
#     source file: a new hyperlinked "Latest Handoff File" cell (column C)
#     is populated, the "Latest Handoff Datetime" (column D) is stamped,
#     and the "Handoff Reason" (column H) becomes "Include".
#   - The second table row (the .localization-config row) is untouched
#     in content, only its relationship id shifts because of the newly
#     inserted hyperlink before it.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/64d8251472525af5cb6fc19cd80e2e38709cb2ee"
$sourceMdUrl = $repoBase + "/e2e/2546b087-f924-469b-bc49-93289baa5b90.md"
$configUrl = $repoBase + "/.localization-config"

function Update-LocalizationSheet {
    param([object]$ws, [string]$handoffFileName, [string]$handoffDateTime)

    # Row 2 status: handoff transform succeeded this time.
    $ws.Range("B2").Value = "Ready for handoff"

    # Re-create the row's hyperlinks in display order (A2, C2, A3) so the
    # new Latest Handoff File link lands between the existing two links,
    # matching how Excel lays out relationship ids when a value is filled
    # in the middle of an already-linked row.
    $ws.Range("A1").Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $sourceMdUrl, "", "", "2546b087-f924-469b-bc49-93289baa5b90.md")

    $handoffUrl = $repoBase + "/e2e/" + $handoffFileName
    $ws.Hyperlinks.Add($ws.Range("C2"), $handoffUrl, "", "", $handoffFileName)

    $ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, "", "", ".localization-config")

    # Latest Handoff Datetime for the newly generated file.
    $ws.Range("D2").Value = $handoffDateTime

    # Handoff Reason: this file is now included in the handoff.
    $ws.Range("H2").Value = "Include"
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LocalizationSheet $wsZhCn `
    "2546b087-f924-469b-bc49-93289baa5b90.db11a40de478fe0108d0d71aebc96cb84d6485f5.zh-cn.xlf" `
    "2016-01-13 16:10:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LocalizationSheet $wsDeDe `
    "2546b087-f924-469b-bc49-93289baa5b90.db11a40de478fe0108d0d71aebc96cb84d6485f5.de-de.xlf" `
    "2016-01-13 16:10:37"
